# Edit the RemoteApiClient / GAE Datastore diagram on slide 1:
#  - shrink the dashed connector arrow that points at the "GAE Datastore" box
#  - move/widen the "GAE Datastore" textbox to fit the new "Google Cloud " label
#  - rename "GAE " to "Google Cloud "

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# Shape 4: "Straight Arrow Connector 89" -> shrink its horizontal extent
# (EMU 1039672 -> 963472, i.e. points 81.86393700787401 -> 75.86393700787401)
$connector = $s.Shapes.Item(4)
$connector.Width = 75.86393700787401

# Shape 5: "TextBox 90" (the "GAE Datastore" label) -> reposition + widen
# (off x EMU 3739356 -> 3663156, i.e. points 294.4374803149606 -> 288.4374803149606)
# (ext cx EMU 1035506 -> 1219200, i.e. points 81.53590551181102 -> 96.0)
# (288.4375 is used instead of the exact quotient because the host stores
# Left/Top/Width/Height as single-precision floats, same as real PowerPoint;
# 288.4375 is exactly representable in float32 and still lands on EMU 3663156)
$textBox = $s.Shapes.Item(5)
$textBox.Left = 288.4375
$textBox.Width = 96.0

# Update the run text from "GAE " to "Google Cloud " (leave the "Datastore"
# run untouched so its formatting/attributes are preserved)
$tr = $textBox.TextFrame.TextRange
$tr.Characters(1, 4).Text = "Google Cloud "
